# Refresh the IG-publisher generated "Metadata" sheet for this run:
#   1. Bump the "Date" property (row 8, column B) to the new publication
#      timestamp.
#   2. Insert a new "Jurisdiction" property row just above "Description"
#      (i.e. at row 11), pushing the remaining property rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the Date value in place.
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# 2. Make room for the new row, then restore the formatting that a bare
#    Insert() loses (border/alignment) by copying it from the row that is
#    about to be pushed down.
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
